# Append a new log row (row 10) to Sheet1, mirroring the SKIPPED rows above.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 10

$ws.Cells.Item($row, 1).Value = "2025-08-14 06:52:19 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-14 12:22:19 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

# Match the style used by the preceding data rows (e.g. row 9) for every cell in the new row.
$srcRange = $ws.Range("A9:H9")
$dstRange = $ws.Range("A10:H10")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 8).Value = ""

$excel.CutCopyMode = $false
